$wb = $excel.ActiveWorkbook

$wsFBS = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# ---- Sheet: FBS ----
$wsFBS.Range("Q2").Value = "WNW"
$wsFBS.Range("Y2").Value = 57.5
$wsFBS.Range("Z2").Value = -115
$wsFBS.Range("AB2").Value = 3.5
$wsFBS.Range("AE2").Value = 0
$wsFBS.Range("AF2").Value = -4.5
$wsFBS.Range("AK2").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK3").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AB4").Value = 8
$wsFBS.Range("AF4").Value = -1.5
$wsFBS.Range("AK4").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Y5").Value = 47.5
$wsFBS.Range("Z5").Value = -115
$wsFBS.Range("AE5").Value = -0.09523809523809523
$wsFBS.Range("AK5").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK6").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK7").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK8").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Y9").Value = 49.5
$wsFBS.Range("AE9").Value = 0
$wsFBS.Range("AK9").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK10").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK11").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Y12").Value = 48.5
$wsFBS.Range("AE12").Value = -0.05825242718446602
$wsFBS.Range("AK12").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK13").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AB14").Value = -7.5
$wsFBS.Range("AF14").Value = 0.5
$wsFBS.Range("AK14").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AB15").Value = 3.5
$wsFBS.Range("AF15").Value = 0
$wsFBS.Range("AK15").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("N16").Value = "NNE"
$wsFBS.Range("O16").Value = 66.77
$wsFBS.Range("Q16").Value = "NNE"
$wsFBS.Range("Y16").Value = 54.5
$wsFBS.Range("AB16").Value = 14.5
$wsFBS.Range("AE16").Value = -0.01801801801801802
$wsFBS.Range("AF16").Value = 1
$wsFBS.Range("AK16").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK17").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Y18").Value = 41.5
$wsFBS.Range("Z18").Value = -118
$wsFBS.Range("AE18").Value = 0.02469135802469136
$wsFBS.Range("AK18").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK19").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AB20").Value = 10
$wsFBS.Range("AF20").Value = -2.5
$wsFBS.Range("AK20").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Y21").Value = 44.5
$wsFBS.Range("AB21").Value = -1
$wsFBS.Range("AE21").Value = -0.04301075268817205
$wsFBS.Range("AF21").Value = 2
$wsFBS.Range("AK21").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK22").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Z23").Value = -110
$wsFBS.Range("AK23").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Q24").Value = "NNW"
$wsFBS.Range("AB24").Value = 3
$wsFBS.Range("AF24").Value = -1
$wsFBS.Range("AK24").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Q25").Value = "NNE"
$wsFBS.Range("Y25").Value = 52.5
$wsFBS.Range("Z25").Value = -112
$wsFBS.Range("AE25").Value = 0.0396039603960396
$wsFBS.Range("AK25").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AB26").Value = -14.5
$wsFBS.Range("AF26").Value = -2.5
$wsFBS.Range("AK26").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK27").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AB28").Value = -13.5
$wsFBS.Range("AF28").Value = -1.5
$wsFBS.Range("AK28").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AB29").Value = -6.5
$wsFBS.Range("AF29").Value = -0.5
$wsFBS.Range("AK29").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Y30").Value = 44.5
$wsFBS.Range("Z30").Value = -110
$wsFBS.Range("AE30").Value = -0.04301075268817205
$wsFBS.Range("AK30").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK31").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Y32").Value = 46.5
$wsFBS.Range("AE32").Value = -0.04123711340206185
$wsFBS.Range("AK32").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Q33").Value = "WSW"
$wsFBS.Range("Z33").Value = -114
$wsFBS.Range("AK33").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK34").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK35").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK36").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Q37").Value = "WNW"
$wsFBS.Range("AB37").Value = 3
$wsFBS.Range("AF37").Value = 0.5
$wsFBS.Range("AK37").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Q38").Value = "WNW"
$wsFBS.Range("AB38").Value = -17
$wsFBS.Range("AF38").Value = 1.5
$wsFBS.Range("AK38").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Y39").Value = 57.5
$wsFBS.Range("Z39").Value = -110
$wsFBS.Range("AE39").Value = -0.04958677685950413
$wsFBS.Range("AK39").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("Y40").Value = 44.5
$wsFBS.Range("Z40").Value = -108
$wsFBS.Range("AE40").Value = 0
$wsFBS.Range("AK40").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK41").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK42").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK43").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK44").Value = "2024-10-31T16:22:59.541428"
$wsFBS.Range("AK45").Value = "2024-10-31T16:22:59.541428"

# ---- Sheet: Other ----
$wsOther.Range("Q15").Value = 74.12
$wsOther.Range("R15").Value = 7.1
$wsOther.Range("S19").Value = "S"
$wsOther.Range("A20").Value = "Lafayette vs Bucknell"
$wsOther.Range("B20").Value = "Bucknell"
$wsOther.Range("C20").Value = "Lafayette"
$wsOther.Range("J20").Value = 69.76284787
$wsOther.Range("K20").Value = 53.28
$wsOther.Range("L20").Value = 53.94
$wsOther.Range("N20").Value = 1924
$wsOther.Range("O20").Value = "SSE"
$wsOther.Range("P20").Value = "SSE"
$wsOther.Range("Q20").Value = 54.79999999999999
$wsOther.Range("S20").Value = "SSE"
$wsOther.Range("X20").Value = "40.9519411, -76.884752"
$wsOther.Range("A21").Value = "William & Mary vs North Carolina A&T"
$wsOther.Range("B21").Value = "North Carolina A&T"
$wsOther.Range("C21").Value = "William & Mary"
$wsOther.Range("J21").Value = 209.04629328
$wsOther.Range("K21").Value = 60.18
$wsOther.Range("L21").Value = 60.45
$wsOther.Range("N21").Value = 1981
$wsOther.Range("O21").Value = "WSW"
$wsOther.Range("P21").Value = "WSW"
$wsOther.Range("Q21").Value = 64.7
$wsOther.Range("R21").Value = 5.4
$wsOther.Range("S21").Value = "WSW"
$wsOther.Range("X21").Value = "36.0814337, -79.7700391"
$wsOther.Range("A22").Value = "Cornell vs Princeton"
$wsOther.Range("B22").Value = "Princeton"
$wsOther.Range("C22").Value = "Cornell"
$wsOther.Range("F22").Value = "High"
$wsOther.Range("J22").Value = -225.31579208
$wsOther.Range("K22").Value = 55.59
$wsOther.Range("L22").Value = 50.04
$wsOther.Range("N22").Value = 1998
$wsOther.Range("O22").Value = "SSE"
$wsOther.Range("P22").Value = "SSE"
$wsOther.Range("Q22").Value = 59.66
$wsOther.Range("R22").Value = 7.6
$wsOther.Range("S22").Value = "SSE"
$wsOther.Range("X22").Value = "40.3457928, -74.6500047"
$wsOther.Range("A23").Value = "Delaware State vs Howard"
$wsOther.Range("B23").Value = "Howard"
$wsOther.Range("C23").Value = "Delaware State"
$wsOther.Range("J23").Value = 42.51443672
$wsOther.Range("K23").Value = 57.45
$wsOther.Range("L23").Value = 57.78
$wsOther.Range("N23").Value = 1926
$wsOther.Range("O23").Value = "S"
$wsOther.Range("P23").Value = "SSE"
$wsOther.Range("Q23").Value = 62.72
$wsOther.Range("R23").Value = 3.5
$wsOther.Range("S23").Value = "SSE"
$wsOther.Range("X23").Value = "38.925487, -77.021047"
$wsOther.Range("A24").Value = "Stony Brook vs Bryant University"
$wsOther.Range("B24").Value = "Bryant University"
$wsOther.Range("C24").Value = "Stony Brook"
$wsOther.Range("J24").Value = 91.83950425
$wsOther.Range("K24").Value = 51.85
$wsOther.Range("L24").Value = 54.66
$wsOther.Range("N24").Value = 1999
$wsOther.Range("O24").Value = "SSE"
$wsOther.Range("Q24").Value = 55.04
$wsOther.Range("R24").Value = 5.6
$wsOther.Range("X24").Value = "41.924206, -71.5385124"
$wsOther.Range("A25").Value = "Norfolk State vs Morgan State"
$wsOther.Range("B25").Value = "Morgan State"
$wsOther.Range("C25").Value = "Norfolk State"
$wsOther.Range("J25").Value = 62.136934276
$wsOther.Range("K25").Value = 56.23
$wsOther.Range("L25").Value = 61.33
$wsOther.Range("N25").Value = 1949
$wsOther.Range("O25").Value = "SSE"
$wsOther.Range("P25").Value = "SSE"
$wsOther.Range("Q25").Value = 61.46
$wsOther.Range("R25").Value = 4.3
$wsOther.Range("S25").Value = "SSE"
$wsOther.Range("X25").Value = "39.3439778, -76.5829573"
$wsOther.Range("A26").Value = "New Hampshire vs Albany"
$wsOther.Range("B26").Value = "Albany"
$wsOther.Range("C26").Value = "New Hampshire"
$wsOther.Range("J26").Value = 45.36820984000001
$wsOther.Range("K26").Value = 50.74
$wsOther.Range("L26").Value = 50.46
$wsOther.Range("N26").Value = 2013
$wsOther.Range("O26").Value = "SE"
$wsOther.Range("P26").Value = "SE"
$wsOther.Range("Q26").Value = 51.68
$wsOther.Range("R26").Value = 6.9
$wsOther.Range("S26").Value = "SE"
$wsOther.Range("X26").Value = "42.6808062, -73.8272841"
$wsOther.Range("A27").Value = "Monmouth vs Rhode Island"
$wsOther.Range("B27").Value = "Rhode Island"
$wsOther.Range("C27").Value = "Monmouth"
$wsOther.Range("J27").Value = 22.78801498
$wsOther.Range("K27").Value = 52.81
$wsOther.Range("L27").Value = 55.89
$wsOther.Range("N27").Value = 1928
$wsOther.Range("O27").Value = "S"
$wsOther.Range("P27").Value = "S"
$wsOther.Range("Q27").Value = 57.2
$wsOther.Range("R27").Value = 4.5
$wsOther.Range("S27").Value = "S"
$wsOther.Range("X27").Value = "41.4879984, -71.5347458"
$wsOther.Range("A28").Value = "Stetson vs Butler"
$wsOther.Range("B28").Value = "Butler"
$wsOther.Range("C28").Value = "Stetson"
$wsOther.Range("F28").Value = "Mid"
$wsOther.Range("J28").Value = 206.54261019
$wsOther.Range("K28").Value = 54.51
$wsOther.Range("L28").Value = 71.79
$wsOther.Range("N28").Value = 1928
$wsOther.Range("O28").Value = "NW"
$wsOther.Range("P28").Value = "NNW"
$wsOther.Range("Q28").Value = 63.26
$wsOther.Range("R28").Value = 4.5
$wsOther.Range("S28").Value = "NNW"
$wsOther.Range("X28").Value = "39.8444694, -86.1659242"
$wsOther.Range("A29").Value = "Dayton vs Presbyterian"
$wsOther.Range("B29").Value = "Presbyterian"
$wsOther.Range("C29").Value = "Dayton"
$wsOther.Range("F29").Value = "High"
$wsOther.Range("J29").Value = -28.2418213
$wsOther.Range("K29").Value = 62.6
$wsOther.Range("L29").Value = 54.68
$wsOther.Range("N29").Value = 2002
$wsOther.Range("O29").Value = "WNW"
$wsOther.Range("P29").Value = "WSW"
$wsOther.Range("Q29").Value = 71.24000000000001
$wsOther.Range("R29").Value = 2.1
$wsOther.Range("S29").Value = "WSW"
$wsOther.Range("X29").Value = "34.46158, -81.861368"
$wsOther.Range("A30").Value = "North Dakota vs Indiana State"
$wsOther.Range("B30").Value = "Indiana State"
$wsOther.Range("C30").Value = "North Dakota"
$wsOther.Range("F30").Value = "Mid"
$wsOther.Range("J30").Value = -99.72589109999998
$wsOther.Range("K30").Value = 55.49
$wsOther.Range("L30").Value = 42.5
$wsOther.Range("N30").Value = 1925
$wsOther.Range("O30").Value = "NW"
$wsOther.Range("P30").Value = "NW"
$wsOther.Range("Q30").Value = 65.3
$wsOther.Range("R30").Value = 6.1
$wsOther.Range("S30").Value = "NW"
$wsOther.Range("X30").Value = "39.4746858, -87.3669599"
$wsOther.Range("A31").Value = "Colgate vs Fordham"
$wsOther.Range("B31").Value = "Fordham"
$wsOther.Range("C31").Value = "Colgate"
$wsOther.Range("J31").ClearContents() | Out-Null
$wsOther.Range("K31").Value = 54.57
$wsOther.Range("L31").Value = 47.85
$wsOther.Range("N31").ClearContents() | Out-Null
$wsOther.Range("O31").Value = "S"
$wsOther.Range("P31").Value = "S"
$wsOther.Range("Q31").Value = 58.63999999999999
$wsOther.Range("R31").Value = 3.9
$wsOther.Range("S31").Value = "S"
$wsOther.Range("X31").Value = "40.8466508, -73.8785937"
$wsOther.Range("A32").Value = "Villanova vs Hampton"
$wsOther.Range("B32").Value = "Hampton"
$wsOther.Range("C32").Value = "Villanova"
$wsOther.Range("J32").Value = -124.809050062
$wsOther.Range("K32").Value = 61.55
$wsOther.Range("L32").Value = 55.05
$wsOther.Range("N32").Value = 1928
$wsOther.Range("O32").Value = "SSW"
$wsOther.Range("P32").Value = "SSW"
$wsOther.Range("Q32").Value = 58.88
$wsOther.Range("R32").Value = 11.4
$wsOther.Range("S32").Value = "SSW"
$wsOther.Range("X32").Value = "37.0200894, -76.3331555"
$wsOther.Range("S38").Value = "NNE"
$wsOther.Range("S45").Value = "NNE"
